# ex9.1.9(Linear) - Strong Stationary generator, alpha non zero
# Update computed values for the non-convex experiment run (all but the 5th).
#
# All these cells hold text (shared-string) representations of numbers /
# algebraic expressions -- not actual numeric cell values -- in the source
# workbook. Setting .Value directly on a numeric-looking string makes Excel
# coerce it to a real number (and bolts on a Text number-format), so instead
# we write a formula that evaluates to the literal text, then demote the
# formula to a plain value via Copy + PasteSpecial(values). That keeps the
# cell type/number format identical to the original (t="s", default style).

function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

$wb = $excel.ActiveWorkbook

# NB: worksheet lookup by name is case-insensitive in this host, and the
# workbook has two sheets whose names differ only by case ("Vector_bf" /
# "Vector_BF") -- Item(name) would collide, so every sheet is addressed by
# its 1-based tab position instead (matches the order in workbook.xml).

# --- Restricciones_del_follower (sheet 3) -------------------------------
$ws3 = $wb.Worksheets.Item(3)

Set-TextValue $ws3.Range("A2") "-2.2665126862411356 - x + 0.8732749362278731y"
Set-TextValue $ws3.Range("B2") "4.2665126862411356"
Set-TextValue $ws3.Range("D2") "0.62"
Set-TextValue $ws3.Range("E2") "5.8999999999999995"
Set-TextValue $ws3.Range("F2") "2.5"

Set-TextValue $ws3.Range("A3") "-16.162969222615605 - 0.25x + 3.388306752564147y"
Set-TextValue $ws3.Range("B3") "14.162969222615605"
Set-TextValue $ws3.Range("D3") "0.96"
Set-TextValue $ws3.Range("E3") "5.0"
Set-TextValue $ws3.Range("F3") "9.7"

Set-TextValue $ws3.Range("A4") "2.552455640104636 + x - 0.931393422785495y"
Set-TextValue $ws3.Range("B4") "-10.552455640104636"
Set-TextValue $ws3.Range("D4") "0.88"
Set-TextValue $ws3.Range("E4") "9.9"
Set-TextValue $ws3.Range("F4") "0.8"

Set-TextValue $ws3.Range("A5") "-19.340469222615607 + x + 3.388306752564147y"
Set-TextValue $ws3.Range("B5") "16.700469222615606"
Set-TextValue $ws3.Range("D5") "0.29"
Set-TextValue $ws3.Range("E5") "8.299999999999999"
Set-TextValue $ws3.Range("F5") "9.7"

Set-TextValue $ws3.Range("A6") "-6.421396745838298 + 1.1527229158207923y"
Set-TextValue $ws3.Range("B6") "5.671396745838298"
Set-TextValue $ws3.Range("D6") "0.34"
Set-TextValue $ws3.Range("E6") "2.5"
Set-TextValue $ws3.Range("F6") "3.3000000000000003"

# --- Punto_modificado (sheet 4) ------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("A2") "2.0300000000000002"
Set-TextValue $ws4.Range("B2") "4.92"

# --- Vector_bf (sheet 5) -----------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-3.3491134804942986"

# --- Vector_BF (sheet 6) -------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "-12.05"
Set-TextValue $ws6.Range("A3") "-44.877814336823185"

# --- Vector_Alpha (sheet 7) ----------------------------------------------
# Unlike the other sheets this one stores an actual numeric cell (no t="s"),
# so a plain numeric assignment is correct here.
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 2.86278684557099
